$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 156, shifting existing rows 156:180 down to 157:181.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new record's data.
$ws.Cells.Item(156, 1).Value = 4
$ws.Cells.Item(156, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(156, 3).Value = "Los Lagos"
$ws.Cells.Item(156, 4).Value = 44491
$ws.Cells.Item(156, 5).Value = 10
$ws.Cells.Item(156, 6).Value = 100112040
$ws.Cells.Item(156, 7).Value = "Cilantro"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 300
$ws.Cells.Item(156, 11).Value = 10000
$ws.Cells.Item(156, 12).Value = 10000
$ws.Cells.Item(156, 13).Value = 10000
$ws.Cells.Item(156, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(156, 15).Value = "Región Metropolitana"
$ws.Cells.Item(156, 16).Value = 278
$ws.Cells.Item(156, 17).Value = 36
$ws.Cells.Item(156, 18).Value = "Hortaliza"
